$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a text value to a cell while preventing Excel from
# auto-coercing numeric-looking strings (e.g. "238.04", "1.000") into
# real numbers. We briefly force Text number format for the write, then
# restore the default "Normal" style so the cell ends up styled exactly
# like the untouched cells around it.
function Set-TextValue($cell, $val) {
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = "30.354.40"
$ws.Range("E2").Value = "  -1.15%  "

$ws.Range("D3").Value = "1.889.38"
$ws.Range("E3").Value = "  -1.48%  "

$ws.Range("E4").Value = "  -0.11%  "

Set-TextValue $ws.Range("D5") "238.04"
$ws.Range("E5").Value = "  -1.12%  "

$ws.Range("E6").Value = "  -0.06%  "

Set-TextValue $ws.Range("D7") "0.4812"
$ws.Range("E7").Value = "  -2.61%  "

Set-TextValue $ws.Range("D8") "0.2896"
$ws.Range("E8").Value = "  -3.73%  "

Set-TextValue $ws.Range("D9") "0.06600"
$ws.Range("E9").Value = "  -2.67%  "

$ws.Range("D10").Value = "1.903.79"
$ws.Range("E10").Value = "  -0.09%  "

Set-TextValue $ws.Range("D11") "16.91"
$ws.Range("E11").Value = "  -1.85%  "

Set-TextValue $ws.Range("D12") "0.07389"
$ws.Range("E12").Value = "  +0.95%  "

Set-TextValue $ws.Range("D13") "5.167"
$ws.Range("E13").Value = "  -1.25%  "

Set-TextValue $ws.Range("D14") "87.81"
$ws.Range("E14").Value = "  -0.80%  "

Set-TextValue $ws.Range("D15") "0.6626"
$ws.Range("E15").Value = "  -1.90%  "

$ws.Range("D16").Value = "30.327.12"

Set-TextValue $ws.Range("D17") "13.45"
$ws.Range("E17").Value = "  -0.86%  "

$ws.Range("B18").Value = "Dai"
$ws.Range("C18").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue $ws.Range("D18") "1.000"
$ws.Range("E18").Value = "  -0.16%  "

$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue $ws.Range("D19") "0.000007749"
$ws.Range("E19").Value = "  -2.94%  "

Set-TextValue $ws.Range("D20") "5.460"
$ws.Range("E20").Value = "  +1.53%  "

$ws.Range("D21").Value = "2.149.14"
$ws.Range("E21").Value = "  -0.74%  "

Set-TextValue $ws.Range("D22") "1.001"
$ws.Range("E22").Value = "  -0.08%  "

Set-TextValue $ws.Range("D23") "192.30"
$ws.Range("E23").Value = "  -1.93%  "

Set-TextValue $ws.Range("D24") "6.211"
$ws.Range("E24").Value = "  -2.07%  "

Set-TextValue $ws.Range("D25") "9.470"
$ws.Range("E25").Value = "  -2.20%  "

Set-TextValue $ws.Range("D26") "165.06"
$ws.Range("E26").Value = "  +1.94%  "

Set-TextValue $ws.Range("D27") "18.22"
$ws.Range("E27").Value = "  -2.27%  "

Set-TextValue $ws.Range("D28") "1.941"
$ws.Range("E28").Value = "  -1.10%  "

$ws.Range("E29").Value = "  -0.94%  "

$ws.Range("E30").Value = "  -2.38%  "

Set-TextValue $ws.Range("D31") "0.09169"
$ws.Range("E31").Value = "  +0.18%  "

Set-TextValue $ws.Range("D32") "4.053"
$ws.Range("E32").Value = "  -0.63%  "

Set-TextValue $ws.Range("D33") "0.05078"
$ws.Range("E33").Value = "  -4.12%  "

Set-TextValue $ws.Range("D34") "0.7317"
$ws.Range("E34").Value = "  -1.74%  "

Set-TextValue $ws.Range("D35") "1.143"
$ws.Range("E35").Value = "  +1.86%  "

Set-TextValue $ws.Range("D36") "2.716"
$ws.Range("E36").Value = "  +0.07%  "

Set-TextValue $ws.Range("D37") "0.01828"
$ws.Range("E37").Value = "  -1.46%  "

$ws.Range("E38").Value = "  -2.82%  "

Set-TextValue $ws.Range("D39") "0.9204"
$ws.Range("E39").Value = "  -1.03%  "

Set-TextValue $ws.Range("D40") "2.091"
$ws.Range("E40").Value = "  +0.16%  "

$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue $ws.Range("D41") "5.911"
$ws.Range("E41").Value = "  -0.50%  "

$ws.Range("B42").Value = "Quant"
$ws.Range("C42").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue $ws.Range("D42") "106.25"
$ws.Range("E42").Value = "  -0.74%  "

Set-TextValue $ws.Range("D43") "0.4334"
$ws.Range("E43").Value = "  -3.80%  "

$ws.Range("E44").Value = "  -0.01%  "

Set-TextValue $ws.Range("D45") "0.1378"
$ws.Range("E45").Value = "  -1.91%  "

Set-TextValue $ws.Range("D46") "7.678"
$ws.Range("E46").Value = "  -0.57%  "

Set-TextValue $ws.Range("D47") "1.573"
$ws.Range("E47").Value = "  +8.25%  "

Set-TextValue $ws.Range("D48") "65.13"
$ws.Range("E48").Value = "  -9.70%  "

Set-TextValue $ws.Range("D49") "8.968"
$ws.Range("E49").Value = "  -1.77%  "

Set-TextValue $ws.Range("D50") "34.26"
$ws.Range("E50").Value = "  -3.39%  "

Set-TextValue $ws.Range("D51") "0.05779"
$ws.Range("E51").Value = "  -2.05%  "
